# internal beads hvb 0.62
# Adds a new data row (row 5) to the "IFCB109" sheet recording the
# 2016-07-23 internal-beads calibration run (hvB = 0.62), and updates the
# window selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IFCB109")
$ws.Activate()

# --- New row 5 values -------------------------------------------------
$ws.Range("A5").Value = 20160723
$ws.Range("B5").Value = 140329
$ws.Range("C5").Value = 0.62
$ws.Range("D5").Value = "na"
$ws.Range("E5").Value = 2.4339
$ws.Range("F5").Value = 0.16643
$ws.Range("G5").Value = "1.5-3.5"
$ws.Range("H5").Value = "0.12-0.25"
$ws.Range("I5").Value = 0.5
$ws.Range("J5").Value = "na"
$ws.Range("K5").Value = 0.081054
$ws.Range("L5").Value = 0.27425
$ws.Range("M5").Value = 0.1538
$ws.Range("N5").Value = 0.24159
$ws.Range("P5").Value = "internal beads use all signals"

# --- Window / selection state ------------------------------------------
# Real Excel scrolls the grid so column H is left-most and selects N6;
# reproduce the reachable part of that (the cell selection) via COM.
$excel.Goto($ws.Range("H1"), $true)
$ws.Range("N6").Select()
